# Updated to include week 14 stats
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - #Autodrafting in 2016 / Emmett
$ws.Range("D2").Value = 42
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 54
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 6
$ws.Range("AA2").Value = 122
$ws.Range("AO2").Value = 112.5

# Row 3 - A Nasty Moses / Max
$ws.Range("D3").Value = 37
$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 59
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 9
$ws.Range("AA3").Value = 70
$ws.Range("AO3").Value = 90.5

# Row 4 - Beauty and the Beast Mode / Eduardo
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 55
$ws.Range("J4").Value = 9
$ws.Range("K4").Value = 5
$ws.Range("AA4").Value = 137
$ws.Range("AO4").Value = 92

# Row 5 - Beejes N Cream / Ranjan
$ws.Range("D5").Value = 44
$ws.Range("E5").Value = 44
$ws.Range("F5").Value = 93
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 8
$ws.Range("AA5").Value = 112
$ws.Range("AO5").Value = 122

# Row 6 - Donte's Winferno / Zach
$ws.Range("D6").Value = 31
$ws.Range("E6").Value = 29
$ws.Range("F6").Value = 82
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 9
$ws.Range("AA6").Value = 113.5
$ws.Range("AO6").Value = 143

# Row 7 - Hyde Yo Kids, Hide Yo Wife / Andy
$ws.Range("D7").Value = 27
$ws.Range("E7").Value = 27
$ws.Range("F7").Value = 58
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 8
$ws.Range("AA7").Value = 78
$ws.Range("AO7").Value = 92.5

# Row 8 - Jizz n' Potatoes / Mike
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = 18
$ws.Range("F8").Value = 51
$ws.Range("J8").Value = 9
$ws.Range("K8").Value = 4
$ws.Range("AA8").Value = 90.5
$ws.Range("AO8").Value = 70

# Row 9 - Mariota Kart 64 / Lauren
$ws.Range("D9").Value = 22
$ws.Range("E9").Value = 21
$ws.Range("F9").Value = 59
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 9
$ws.Range("AA9").Value = 92
$ws.Range("AO9").Value = 137

# Row 10 - Mdphd Killa!!! / Jeremy
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 14
$ws.Range("F10").Value = 44
$ws.Range("J10").Value = 9
$ws.Range("K10").Value = 5
$ws.Range("AA10").Value = 92.5
$ws.Range("AO10").Value = 78

# Row 11 - PASSWORD IS TACO / Ted
$ws.Range("D11").Value = 33
$ws.Range("E11").Value = 32
$ws.Range("F11").Value = 74
$ws.Range("J11").Value = 8
$ws.Range("K11").Value = 6
$ws.Range("AA11").Value = 143
$ws.Range("AO11").Value = 113.5

# Update the view so the active/selected cell reflects the last edit
$null = $ws.Range("AO11").Select()
